$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oktober")

# --- Fill in hours + details for rows 24-26 and 28 ---
$ws.Range("B24").Value = 2.5
$ws.Range("C24").Value = "CC2650 in betrieb genommen, TI RTOS eingelesen, "

$ws.Range("B25").Value = 2
$ws.Range("C25").Value = "I2C Beispiel geschrieben, gibt fehler beim start der Transaction "

$ws.Range("B26").Value = 1.5
$ws.Range("C26").Value = "CC2650 I2C Beispiel laufen lassen, nun spinnt der debugger "

$ws.Range("B28").Value = 1.5
$ws.Range("C28").Value = "I/O Pins geprüft zwischen CC2650 und Senor Hub, Bild erstellt, wie Pins angepasst werden müssen. Konkret I2C  muss anders liegen, wenn Analog kein interrupt hat muss der Lichtsensor Pin umgelegt werden.  Stecker gesucht und flachbandkabel! "

# C28 wraps across several lines like C5/C16, so enable wrap text + taller row
$ws.Range("C28").WrapText = $true
$ws.Rows.Item(28).RowHeight = 63

# Rows 5 and 16 previously had extra height to fit wrapped text; now auto-fit back down
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(16).AutoFit()

# --- Update the view: scroll position and active selection ---
$ws.Application.GoTo($ws.Range("A9")) | Out-Null
$ws.Range("E24").Select() | Out-Null
